$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Value = "Luca Giordani"
$ws.Range("B19").Value = "Alberto Cerisara | Shark Attack"
$ws.Range("C19").Value = "Mattia Baldessarini | Shark Attack"
$ws.Range("D19").Value = "ALESSIO FARINATI | Pinguini Trentini"
$ws.Range("E19").Value = "Antonio Calabrò | Avanzi"
$ws.Range("F19").Value = "Leonardo Fedrigotti | Nazionale Ferrovieri"
